$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet currently holds rows 1 (header) + 2-41 (schedule for group "B2-B").
# We are appending 40 more schedule rows (42-81) for group "B2-C", using the
# exact same alternating row style as the existing data: the first row of
# each pair uses styles 2/2/2/2/3/4/5 (cols A-G) and the second row of the
# pair uses styles 6/6/6/6/7/8/9 - i.e. rows 40 and 41 are a perfect 2-row
# template for the whole new block.
# ---------------------------------------------------------------------------

# 1) Stamp that alternating style pattern across rows 42-81 in one shot by
#    tiling the existing formatted 2-row template (rows 40:41). This copies
#    both values and formats, so every new cell immediately has the correct
#    style index; the (template) values get overwritten below.
$ws.Range("A40:G41").Copy($ws.Range("A42:G81"))

# 2) Scratch cell (far outside the used range, cleared at the end) that we
#    pre-format as Text ("@"). Excel's normal smart-typing would otherwise
#    parse strings like "03/01/2026" or "1" into a date serial / number as
#    soon as they are assigned to a cell - exactly as it does interactively.
#    Pre-formatting this helper cell as Text keeps the assigned string
#    literal, and copying *only its value* (xlPasteValues = -4163) into the
#    destination cell avoids carrying the helper's own Text style along.
#    We then paste-special *only the format* (xlPasteFormats = -4122) from
#    the matching template cell (row 40 or 41) onto the destination so it
#    keeps the same style index as the rest of that row-pair.
$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"

function Set-LiteralText($targetAddr, $templateAddr, $text) {
    $scratch.Value = $text
    $scratch.Copy()
    $ws.Range($targetAddr).PasteSpecial(-4163)
    $ws.Range($templateAddr).Copy()
    $ws.Range($targetAddr).PasteSpecial(-4122)
}

$newRows = @(
    @{Row=42; Group="B2-C"; Subject="endocrinology"; Session="1"; Date="03/01/2026"}
    @{Row=43; Group="B2-C"; Subject="endocrinology"; Session="2"; Date="04/01/2026"}
    @{Row=44; Group="B2-C"; Subject="endocrinology"; Session="3"; Date="05/01/2026"}
    @{Row=45; Group="B2-C"; Subject="endocrinology"; Session="4"; Date="06/01/2026"}
    @{Row=46; Group="B2-C"; Subject="endocrinology"; Session="5"; Date="07/01/2026"}
    @{Row=47; Group="B2-C"; Subject="endocrinology"; Session="6"; Date="10/01/2026"}
    @{Row=48; Group="B2-C"; Subject="endocrinology"; Session="7"; Date="11/01/2026"}
    @{Row=49; Group="B2-C"; Subject="endocrinology"; Session="8"; Date="12/01/2026"}
    @{Row=50; Group="B2-C"; Subject="endocrinology"; Session="9"; Date="13/01/2026"}
    @{Row=51; Group="B2-C"; Subject="endocrinology"; Session="10"; Date="14/01/2026"}
    @{Row=52; Group="B2-C"; Subject="gastroenterology"; Session="1"; Date="17/01/2026"}
    @{Row=53; Group="B2-C"; Subject="gastroenterology"; Session="2"; Date="18/01/2026"}
    @{Row=54; Group="B2-C"; Subject="gastroenterology"; Session="3"; Date="19/01/2026"}
    @{Row=55; Group="B2-C"; Subject="gastroenterology"; Session="4"; Date="20/01/2026"}
    @{Row=56; Group="B2-C"; Subject="gastroenterology"; Session="5"; Date="21/01/2026"}
    @{Row=57; Group="B2-C"; Subject="gastroenterology"; Session="6"; Date="07/02/2026"}
    @{Row=58; Group="B2-C"; Subject="gastroenterology"; Session="7"; Date="08/02/2026"}
    @{Row=59; Group="B2-C"; Subject="gastroenterology"; Session="8"; Date="09/02/2026"}
    @{Row=60; Group="B2-C"; Subject="gastroenterology"; Session="9"; Date="10/02/2026"}
    @{Row=61; Group="B2-C"; Subject="gastroenterology"; Session="10"; Date="11/02/2026"}
    @{Row=62; Group="B2-C"; Subject="neurology"; Session="1"; Date="20/12/2025"}
    @{Row=63; Group="B2-C"; Subject="neurology"; Session="2"; Date="21/12/2025"}
    @{Row=64; Group="B2-C"; Subject="neurology"; Session="3"; Date="22/12/2025"}
    @{Row=65; Group="B2-C"; Subject="neurology"; Session="4"; Date="23/12/2025"}
    @{Row=66; Group="B2-C"; Subject="neurology"; Session="5"; Date="27/12/2025"}
    @{Row=67; Group="B2-C"; Subject="neurology"; Session="6"; Date="28/12/2025"}
    @{Row=68; Group="B2-C"; Subject="neurology"; Session="7"; Date="29/12/2025"}
    @{Row=69; Group="B2-C"; Subject="neurology"; Session="8"; Date="30/12/2025"}
    @{Row=70; Group="B2-C"; Subject="physical medicine"; Session="1"; Date="24/12/2025"}
    @{Row=71; Group="B2-C"; Subject="physical medicine"; Session="2"; Date="31/12/2025"}
    @{Row=72; Group="B2-C"; Subject="rheumatology"; Session="1"; Date="06/12/2025"}
    @{Row=73; Group="B2-C"; Subject="rheumatology"; Session="2"; Date="07/12/2025"}
    @{Row=74; Group="B2-C"; Subject="rheumatology"; Session="3"; Date="08/12/2025"}
    @{Row=75; Group="B2-C"; Subject="rheumatology"; Session="4"; Date="09/12/2025"}
    @{Row=76; Group="B2-C"; Subject="rheumatology"; Session="5"; Date="10/12/2025"}
    @{Row=77; Group="B2-C"; Subject="rheumatology"; Session="6"; Date="13/12/2025"}
    @{Row=78; Group="B2-C"; Subject="rheumatology"; Session="7"; Date="14/12/2025"}
    @{Row=79; Group="B2-C"; Subject="rheumatology"; Session="8"; Date="15/12/2025"}
    @{Row=80; Group="B2-C"; Subject="rheumatology"; Session="9"; Date="16/12/2025"}
    @{Row=81; Group="B2-C"; Subject="rheumatology"; Session="10"; Date="17/12/2025"}
)

foreach ($row in $newRows) {
    $r = $row.Row
    # Template row alternates 40/41 in lock-step with the existing 2-row
    # style pattern (even destination row -> style-pattern of row 40, odd -> 41).
    if (($r % 2) -eq 0) { $tmpl = 40 } else { $tmpl = 41 }

    $ws.Range("A$r").Value = "Year 5"
    $ws.Range("B$r").Value = $row.Group
    $ws.Range("C$r").Value = $row.Subject

    Set-LiteralText "D$r" "D$tmpl" $row.Session
    Set-LiteralText "E$r" "E$tmpl" $row.Date
    Set-LiteralText "F$r" "F$tmpl" "09:00:00"

    $ws.Range("G$r").Value = 360
}

# 3) Remove the scratch cell's content + formatting entirely so it does not
#    widen the sheet's used range / dimension beyond G81.
$scratch.Clear()
